$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing F/G column values per diff (rows 306-550)
$ws.Cells.Item(306, 6).Value = 74487
$ws.Cells.Item(310, 6).Value = 76476
$ws.Cells.Item(386, 6).Value = 183536
$ws.Cells.Item(387, 6).Value = 351959
$ws.Cells.Item(398, 6).Value = 300863
$ws.Cells.Item(400, 6).Value = 150064
$ws.Cells.Item(407, 6).Value = 158575
$ws.Cells.Item(434, 6).Value = 79512
$ws.Cells.Item(440, 6).Value = 73944
$ws.Cells.Item(456, 6).Value = 50485
$ws.Cells.Item(477, 6).Value = 37219
$ws.Cells.Item(485, 6).Value = 14132
$ws.Cells.Item(487, 6).Value = 6887
$ws.Cells.Item(488, 6).Value = 6405
$ws.Cells.Item(489, 6).Value = 13003
$ws.Cells.Item(490, 6).Value = 10881
$ws.Cells.Item(491, 6).Value = 9959
$ws.Cells.Item(492, 6).Value = 14425
$ws.Cells.Item(493, 6).Value = 8353
$ws.Cells.Item(494, 6).Value = 6785
$ws.Cells.Item(494, 7).Value = 8
$ws.Cells.Item(495, 6).Value = 10495
$ws.Cells.Item(496, 6).Value = 8311
$ws.Cells.Item(497, 6).Value = 7843
$ws.Cells.Item(498, 6).Value = 9245
$ws.Cells.Item(499, 6).Value = 11591
$ws.Cells.Item(500, 6).Value = 7851
$ws.Cells.Item(501, 6).Value = 5845
$ws.Cells.Item(502, 6).Value = 10703
$ws.Cells.Item(503, 6).Value = 7652
$ws.Cells.Item(504, 6).Value = 7621
$ws.Cells.Item(505, 6).Value = 8717
$ws.Cells.Item(506, 6).Value = 11107
$ws.Cells.Item(507, 6).Value = 7414
$ws.Cells.Item(508, 6).Value = 5860
$ws.Cells.Item(509, 6).Value = 9797
$ws.Cells.Item(509, 7).Value = 22
$ws.Cells.Item(510, 6).Value = 7995
$ws.Cells.Item(512, 6).Value = 8672
$ws.Cells.Item(513, 6).Value = 10563
$ws.Cells.Item(514, 6).Value = 7142
$ws.Cells.Item(516, 6).Value = 9461
$ws.Cells.Item(520, 6).Value = 10339
$ws.Cells.Item(523, 6).Value = 10219
$ws.Cells.Item(524, 6).Value = 7838
$ws.Cells.Item(525, 6).Value = 7626
$ws.Cells.Item(526, 6).Value = 8784
$ws.Cells.Item(527, 6).Value = 11492
$ws.Cells.Item(528, 6).Value = 8034
$ws.Cells.Item(529, 6).Value = 5673
$ws.Cells.Item(530, 6).Value = 12690
$ws.Cells.Item(531, 6).Value = 9193
$ws.Cells.Item(532, 6).Value = 10200
$ws.Cells.Item(533, 6).Value = 11809
$ws.Cells.Item(534, 6).Value = 16661
$ws.Cells.Item(535, 6).Value = 10066
$ws.Cells.Item(536, 6).Value = 7913
$ws.Cells.Item(537, 6).Value = 13536
$ws.Cells.Item(538, 6).Value = 11143
$ws.Cells.Item(539, 6).Value = 10433
$ws.Cells.Item(540, 6).Value = 12354
$ws.Cells.Item(541, 6).Value = 16373
$ws.Cells.Item(542, 6).Value = 10190
$ws.Cells.Item(543, 6).Value = 4611
$ws.Cells.Item(544, 6).Value = 14184
$ws.Cells.Item(545, 6).Value = 16506
$ws.Cells.Item(546, 6).Value = 3730
$ws.Cells.Item(547, 6).Value = 13766
$ws.Cells.Item(547, 7).Value = 150
$ws.Cells.Item(548, 6).Value = 16616
$ws.Cells.Item(548, 7).Value = 146
$ws.Cells.Item(549, 6).Value = 9977
$ws.Cells.Item(549, 7).Value = 71
$ws.Cells.Item(550, 6).Value = 8049
$ws.Cells.Item(550, 7).Value = 85

# Append new rows 551-553 (2021-09-06 .. 2021-09-08 data)
$ws.Cells.Item(551, 1).Value = 44445
$ws.Cells.Item(551, 2).Value = 396487
$ws.Cells.Item(551, 3).Value = 8113
$ws.Cells.Item(551, 4).Value = 306
$ws.Cells.Item(551, 5).Value = 12552
$ws.Cells.Item(551, 6).Value = 17047
$ws.Cells.Item(551, 7).Value = 192
$ws.Cells.Item(552, 1).Value = 44446
$ws.Cells.Item(552, 2).Value = 396904
$ws.Cells.Item(552, 3).Value = 9450
$ws.Cells.Item(552, 4).Value = 417
$ws.Cells.Item(552, 5).Value = 12553
$ws.Cells.Item(552, 6).Value = 14680
$ws.Cells.Item(552, 7).Value = 169
$ws.Cells.Item(553, 1).Value = 44447
$ws.Cells.Item(553, 2).Value = 397382
$ws.Cells.Item(553, 3).Value = 34226
$ws.Cells.Item(553, 4).Value = 478
$ws.Cells.Item(553, 5).Value = 12556
$ws.Cells.Item(553, 6).Value = 10404
$ws.Cells.Item(553, 7).Value = 112
